$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 120, shifting existing rows 120:125 down to 121:126.
$ws.Rows(120).Insert()

# Populate the newly inserted row 120 with the new price-report entry.
$ws.Cells.Item(120, 1).Value = 7
$ws.Cells.Item(120, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(120, 3).Value = "Ñuble"
$ws.Cells.Item(120, 4).Value = 45021
$ws.Cells.Item(120, 5).Value = 16
$ws.Cells.Item(120, 6).Value = 100112037
$ws.Cells.Item(120, 7).Value = "Cebollín"
$ws.Cells.Item(120, 8).Value = "Sin especificar"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 120
$ws.Cells.Item(120, 11).Value = 7000
$ws.Cells.Item(120, 12).Value = 7000
$ws.Cells.Item(120, 13).Value = 7000
$ws.Cells.Item(120, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(120, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(120, 16).Value = 194
$ws.Cells.Item(120, 17).Value = 36
$ws.Cells.Item(120, 18).Value = "Hortaliza"
